# Commit: "Added 'Stored in refrigerator' as an option for 'preparation condition'"
#
# This workbook is a CEDAR-style metadata template. The "preparation_condition"
# sheet backs the data-validation dropdown used by the "preparation_condition"
# column (H) on the "Sample Section" sheet. We add the new allowed value as a
# new row at the bottom of that lookup sheet, widen the dropdown's validation
# range to include it, and bump the template's recorded "createdOn" timestamp
# (tracked on the ".metadata" sheet) the same way the template generator does
# whenever the schema changes.

$wb = $excel.ActiveWorkbook

# 1) Add the new "Stored in refrigerator" option to the preparation_condition
#    lookup sheet (it previously had 7 options in rows 1-7; the new one goes
#    into row 8).
$wsCondition = $wb.Sheets.Item("preparation_condition")
$wsCondition.Range("A8").Value = "Stored in refrigerator"
$wsCondition.Range("B8").Value = "https://purl.humanatlas.io/vocab/hravs#HRAVS_0000104"

# 2) Extend the "preparation_condition" column's dropdown validation on the
#    main "Sample Section" sheet so it covers the new row (A1:A7 -> A1:A8).
$wsMain = $wb.Sheets.Item("Sample Section")
$conditionRange = $wsMain.Range("H2:H1001")
$conditionRange.Validation.Modify(3, 1, 1, "'preparation_condition'!`$A`$1:`$A`$8")

# 3) The template's ".metadata" sheet records when the template was generated
#    (pav:createdOn). Update it to reflect this regeneration.
$wsMeta = $wb.Sheets.Item(".metadata")
$wsMeta.Range("C2").Value = "2024-02-21T09:25:13-08:00"
